$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers
$ws.Range("A1").Value = "Week"
$ws.Range("B1").Value = "Predicted_Quantity"

# New data: Week label, Predicted_Quantity
$weeks = @("2025-W43", "2025-W44", "2025-W45", "2025-W46", "2025-W47", "2025-W48", "2025-W49", "2025-W50")
$values = @(21, 16, 12, 21, 19, 20, 28, 25)

for ($i = 0; $i -lt $weeks.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $weeks[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Remove old column C entirely (shift the range so dimension becomes A1:B9)
$ws.Range("C1:C9").Delete()
